$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9743.182000000001
$ws.Range("I62").Value = 7699.2
$ws.Range("J62").Value = 11446.5
$ws.Range("K62").Value = 7699.2
$ws.Range("L62").Value = 11446.5
$ws.Range("M62").Value = -7075.2
$ws.Range("N62").Value = -12694.5

$ws.Range("H65").Value = 9743.182000000001
$ws.Range("I65").Value = 7699.2
$ws.Range("J65").Value = 11446.5
$ws.Range("K65").Value = 38496
$ws.Range("L65").Value = 57232.5
$ws.Range("M65").Value = -35376
$ws.Range("N65").Value = -63472.5

$ws.Range("H86").Value = 153848530
$ws.Range("I86").Value = 125002560
$ws.Range("K86").Value = 125002560
$ws.Range("M86").Value = -125001437

$ws.Range("H89").Value = 153848530
$ws.Range("I89").Value = 125002560
$ws.Range("K89").Value = 625012800
$ws.Range("M89").Value = -625007184

$ws.Range("H106").Value = 11112961
$ws.Range("I106").Value = 11112961
$ws.Range("K106").Value = 11112961
$ws.Range("M106").Value = -11112330

$ws.Range("H132").Value = 3490.3618
$ws.Range("I132").Value = 3458.5476
$ws.Range("K132").Value = 10375.6428
$ws.Range("M132").Value = -7845.6428

$ws.Range("H137").Value = 2198.5
$ws.Range("I137").Value = 867.46155
$ws.Range("K137").Value = 2602.38465
$ws.Range("M137").Value = -52.38464999999997

$ws.Range("H138").Value = 2918.224
$ws.Range("I138").Value = 2884.1428
$ws.Range("J138").Value = 2929.068
$ws.Range("K138").Value = 8652.428400000001
$ws.Range("L138").Value = 8787.204000000002
$ws.Range("M138").Value = -3512.428400000001
$ws.Range("N138").Value = -19067.204

$ws.Range("H141").Value = 3620.625
$ws.Range("J141").Value = 5000
$ws.Range("L141").Value = 15000
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 500945
$ws.Range("I8").Value = 1890
$ws.Range("K8").Value = 1890
$ws.Range("M8").Value = -1746

$ws.Range("H11").Value = 750
$ws.Range("I11").Value = 750
$ws.Range("K11").Value = 750
$ws.Range("M11").Value = -606

$ws.Range("H13").Value = 80000000
$ws.Range("I13").Value = 80000000
$ws.Range("K13").Value = 80000000
$ws.Range("M13").Value = -79999856

$ws.Range("H29").Value = 14501.1875
$ws.Range("I29").Value = 8009
$ws.Range("K29").Value = 8009
$ws.Range("M29").Value = -7701

$ws.Range("H32").Value = 10419316
$ws.Range("I32").Value = 5748714.5
$ws.Range("J32").Value = 22732720
$ws.Range("K32").Value = 5748714.5
$ws.Range("L32").Value = 22732720
$ws.Range("M32").Value = -5748427.5
$ws.Range("N32").Value = -22733294

$ws.Range("H33").Value = 26513
$ws.Range("I33").Value = 20026
$ws.Range("J33").Value = 33000
$ws.Range("K33").Value = 20026
$ws.Range("L33").Value = 33000
$ws.Range("M33").Value = -19697
$ws.Range("N33").Value = -33658

$ws.Range("H56").Value = 50110
$ws.Range("J56").Value = 50110
$ws.Range("L56").Value = 50110
$ws.Range("N56").Value = -51594

$ws.Range("H97").Value = 1230.0322
$ws.Range("J97").Value = 1556.909
$ws.Range("L97").Value = 1556.909
$ws.Range("N97").Value = -2548.909

$ws.Range("H102").Value = 793.94446
$ws.Range("I102").Value = 608.7692
$ws.Range("K102").Value = 608.7692
$ws.Range("M102").Value = 1013.2308

$ws.Range("H122").Value = 3529.926
$ws.Range("I122").Value = 2696.7585
$ws.Range("J122").Value = 5630.9565
$ws.Range("K122").Value = 8090.2755
$ws.Range("L122").Value = 16892.8695
$ws.Range("M122").Value = -5640.2755
$ws.Range("N122").Value = -21792.8695

$ws.Range("H138").Value = 99194.5
$ws.Range("J138").Value = 99194.5
$ws.Range("L138").Value = 99194.5
$ws.Range("N138").Value = -109474.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H86").Value = 2657.4375
$ws.Range("I86").Value = 3621.3333
$ws.Range("J86").Value = 2079.1
$ws.Range("K86").Value = 3621.3333
$ws.Range("L86").Value = 2079.1
$ws.Range("M86").Value = -2498.3333
$ws.Range("N86").Value = -4325.1

$ws.Range("H89").Value = 2657.4375
$ws.Range("I89").Value = 3621.3333
$ws.Range("J89").Value = 2079.1
$ws.Range("K89").Value = 18106.6665
$ws.Range("L89").Value = 10395.5
$ws.Range("M89").Value = -12490.6665
$ws.Range("N89").Value = -21627.5

$ws.Range("H94").Value = 1496.7858
$ws.Range("I94").Value = 897
$ws.Range("J94").Value = 2576.4
$ws.Range("K94").Value = 897
$ws.Range("L94").Value = 2576.4
$ws.Range("M94").Value = -446
$ws.Range("N94").Value = -3478.4

$ws.Range("H105").Value = 3771.7693
$ws.Range("I105").Value = 2428
$ws.Range("J105").Value = 4174.9
$ws.Range("K105").Value = 2428
$ws.Range("L105").Value = 4174.9
$ws.Range("M105").Value = -681
$ws.Range("N105").Value = -7668.9

$ws.Range("H107").Value = 1397.2727
$ws.Range("I107").Value = 1095.2222
$ws.Range("K107").Value = 1095.2222
$ws.Range("M107").Value = 824.7778000000001

$ws.Range("H134").Value = 10941285
$ws.Range("I134").Value = 2166188.2
$ws.Range("K134").Value = 6498564.600000001
$ws.Range("M134").Value = -6496029.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3100.8215
$ws.Range("J58").Value = 4492.625
$ws.Range("L58").Value = 4492.625
$ws.Range("N58").Value = -4898.625

$ws.Range("H98").Value = 73150.60000000001
$ws.Range("J98").Value = 73150.60000000001
$ws.Range("L98").Value = 73150.60000000001
$ws.Range("N98").Value = -77642.60000000001

$ws.Range("H136").Value = 3100.8215
$ws.Range("J136").Value = 4492.625
$ws.Range("L136").Value = 13477.875
$ws.Range("N136").Value = -18577.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2274.4285
$ws.Range("I80").Value = 850.5
$ws.Range("K80").Value = 850.5
$ws.Range("M80").Value = 147.5

$ws.Range("H83").Value = 2274.4285
$ws.Range("I83").Value = 850.5
$ws.Range("K83").Value = 4252.5
$ws.Range("M83").Value = 739.5

$ws.Range("H97").Value = 2649.75
$ws.Range("I97").Value = 1100
$ws.Range("J97").Value = 3166.3333
$ws.Range("K97").Value = 1100
$ws.Range("L97").Value = 3166.3333
$ws.Range("N97").Value = -4158.3333
$ws.Range("M97").Value = -604

$ws.Range("H136").Value = 32788.6
$ws.Range("J136").Value = 32788.6
$ws.Range("L136").Value = 98365.79999999999
$ws.Range("N136").Value = -103465.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 8997.5
$ws.Range("J34").Value = 8663.333000000001
$ws.Range("L34").Value = 8663.333000000001
$ws.Range("N34").Value = -9007.333000000001

$ws.Range("H40").Value = 7160.25
$ws.Range("I40").Value = 7022.6113
$ws.Range("K40").Value = 7022.6113
$ws.Range("M40").Value = -6886.6113

$ws.Range("H82").Value = 3375.2354
$ws.Range("I82").Value = 2916.4
$ws.Range("K82").Value = 2916.4
$ws.Range("M82").Value = -2555.4

$ws.Range("H85").Value = 3375.2354
$ws.Range("I85").Value = 2916.4
$ws.Range("K85").Value = 2916.4
$ws.Range("M85").Value = -1668.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 20015
$ws.Range("I21").Value = 20015
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 20015
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -19780
$ws.Range("N21").ClearContents()

$ws.Range("H24").Value = 30009
$ws.Range("I24").Value = 30009
$ws.Range("K24").Value = 30009
$ws.Range("M24").Value = -29779

$ws.Range("H28").Value = 3771070.2
$ws.Range("I28").Value = 10006341
$ws.Range("J28").Value = 29907.6
$ws.Range("K28").Value = 10006341
$ws.Range("L28").Value = 29907.6
$ws.Range("M28").Value = -10005993
$ws.Range("N28").Value = -30603.6

$ws.Range("H35").Value = 20015
$ws.Range("I35").Value = 20015
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 20015
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -19725
$ws.Range("N35").ClearContents()

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H62").Value = 4736.231
$ws.Range("I62").Value = 3846.375
$ws.Range("J62").Value = 6160
$ws.Range("K62").Value = 3846.375
$ws.Range("L62").Value = 6160
$ws.Range("M62").Value = -3222.375
$ws.Range("N62").Value = -7408

$ws.Range("H65").Value = 4736.231
$ws.Range("I65").Value = 3846.375
$ws.Range("J65").Value = 6160
$ws.Range("K65").Value = 19231.875
$ws.Range("L65").Value = 30800
$ws.Range("M65").Value = -16111.875
$ws.Range("N65").Value = -37040

$ws.Range("H122").Value = 2520.9167
$ws.Range("I122").Value = 2113.025
$ws.Range("K122").Value = 6339.075000000001
$ws.Range("M122").Value = -3889.075000000001

$ws.Range("H132").Value = 2018.2778
$ws.Range("I132").Value = 1532.0851
$ws.Range("K132").Value = 4596.2553
$ws.Range("M132").Value = -2066.2553

$ws.Range("H136").Value = 2560.585
$ws.Range("I136").Value = 1657.5
$ws.Range("J136").Value = 4848.4
$ws.Range("K136").Value = 4972.5
$ws.Range("L136").Value = 14545.2
$ws.Range("M136").Value = -2422.5
$ws.Range("N136").Value = -19645.2
